$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = ""
$ws.Range("N27").Value = ""

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = ""
$ws.Range("N51").Value = ""

$ws.Range("H64").Value = 3167.4443
$ws.Range("I64").Value = 2919.8
$ws.Range("J64").Value = 3262.6924
$ws.Range("K64").Value = 2919.8
$ws.Range("L64").Value = 3262.6924
$ws.Range("M64").Value = -2671.8
$ws.Range("N64").Value = -3758.6924

$ws.Range("H67").Value = 3167.4443
$ws.Range("I67").Value = 2919.8
$ws.Range("J67").Value = 3262.6924
$ws.Range("K67").Value = 2919.8
$ws.Range("L67").Value = 3262.6924
$ws.Range("M67").Value = -2061.8
$ws.Range("N67").Value = -4978.6924

$ws.Range("H86").Value = 6152.55
$ws.Range("I86").Value = 1239.1538
$ws.Range("J86").Value = 15277.429
$ws.Range("K86").Value = 1239.1538
$ws.Range("L86").Value = 15277.429
$ws.Range("M86").Value = -116.1538
$ws.Range("N86").Value = -17523.429

$ws.Range("H89").Value = 6152.55
$ws.Range("I89").Value = 1239.1538
$ws.Range("J89").Value = 15277.429
$ws.Range("K89").Value = 6195.769
$ws.Range("L89").Value = 76387.145
$ws.Range("M89").Value = -579.7690000000002
$ws.Range("N89").Value = -87619.145

$ws.Range("H103").Value = 417049.66
$ws.Range("I103").Value = 833669.3
$ws.Range("J103").Value = 430
$ws.Range("K103").Value = 2501007.9
$ws.Range("L103").Value = 1290
$ws.Range("M103").Value = -2500421.9
$ws.Range("N103").Value = -2462

$ws.Range("H106").Value = 9526656
$ws.Range("I106").Value = 22224066
$ws.Range("J106").Value = 3598.05
$ws.Range("K106").Value = 22224066
$ws.Range("L106").Value = 3598.05
$ws.Range("M106").Value = -22223435
$ws.Range("N106").Value = -4860.05

$ws.Range("H132").Value = 2731.658
$ws.Range("I132").Value = 3099.9033
$ws.Range("J132").Value = 1100.8572
$ws.Range("K132").Value = 9299.7099
$ws.Range("L132").Value = 3302.5716
$ws.Range("M132").Value = -6769.7099
$ws.Range("N132").Value = -8362.571599999999

$ws.Range("H139").Value = 50148
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 50148
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 50148
$ws.Range("M139").Value = ""
$ws.Range("N139").Value = -60428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5253
$ws.Range("I32").Value = 4587.6816
$ws.Range("J32").Value = 8912.25
$ws.Range("K32").Value = 4587.6816
$ws.Range("L32").Value = 8912.25
$ws.Range("M32").Value = -4300.6816
$ws.Range("N32").Value = -9486.25

$ws.Range("H45").Value = 2319.524
$ws.Range("I45").Value = 1666.5333
$ws.Range("J45").Value = 3952
$ws.Range("K45").Value = 1666.5333
$ws.Range("L45").Value = 3952
$ws.Range("M45").Value = -1289.5333
$ws.Range("N45").Value = -4706

$ws.Range("H61").Value = 1485.4706
$ws.Range("I61").Value = 1456.3572
$ws.Range("J61").Value = 1621.3334
$ws.Range("K61").Value = 1456.3572
$ws.Range("L61").Value = 1621.3334
$ws.Range("M61").Value = -1244.3572
$ws.Range("N61").Value = -2045.3334

$ws.Range("H63").Value = 2012.3334
$ws.Range("I63").Value = 2031.7273
$ws.Range("J63").Value = 1799
$ws.Range("K63").Value = 2031.7273
$ws.Range("L63").Value = 1799
$ws.Range("M63").Value = -1345.7273
$ws.Range("N63").Value = -3171

$ws.Range("H66").Value = 2012.3334
$ws.Range("I66").Value = 2031.7273
$ws.Range("J66").Value = 1799
$ws.Range("K66").Value = 10158.6365
$ws.Range("L66").Value = 8995
$ws.Range("M66").Value = -6726.636500000001
$ws.Range("N66").Value = -15859

$ws.Range("H102").Value = 1181.7
$ws.Range("I102").Value = 1127.125
$ws.Range("J102").Value = 1400
$ws.Range("K102").Value = 1127.125
$ws.Range("L102").Value = 1400
$ws.Range("M102").Value = 494.875
$ws.Range("N102").Value = -4644

$ws.Range("H110").Value = 657.1667
$ws.Range("I110").Value = 611.125
$ws.Range("J110").Value = 749.25
$ws.Range("K110").Value = 611.125
$ws.Range("L110").Value = 749.25
$ws.Range("M110").Value = 1433.875
$ws.Range("N110").Value = -4839.25

$ws.Range("H136").Value = 1485.4706
$ws.Range("I136").Value = 1456.3572
$ws.Range("J136").Value = 1621.3334
$ws.Range("K136").Value = 4369.071599999999
$ws.Range("L136").Value = 4864.0002
$ws.Range("M136").Value = -1819.071599999999
$ws.Range("N136").Value = -9964.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1740
$ws.Range("I99").Value = 1800
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 1800
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -302
$ws.Range("N99").Value = -4496

$ws.Range("H105").Value = 1830
$ws.Range("I105").Value = 1622.3529
$ws.Range("J105").Value = 2015.7894
$ws.Range("K105").Value = 1622.3529
$ws.Range("L105").Value = 2015.7894
$ws.Range("M105").Value = 124.6470999999999
$ws.Range("N105").Value = -5509.7894

$ws.Range("H134").Value = 3010.6956
$ws.Range("I134").Value = 2945.0476
$ws.Range("J134").Value = 3700
$ws.Range("K134").Value = 8835.1428
$ws.Range("L134").Value = 11100
$ws.Range("M134").Value = -6300.1428
$ws.Range("N134").Value = -16170

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3485.0244
$ws.Range("I31").Value = 3427.7273
$ws.Range("J31").Value = 3506.0334
$ws.Range("K31").Value = 3427.7273
$ws.Range("L31").Value = 3506.0334
$ws.Range("M31").Value = -3132.7273
$ws.Range("N31").Value = -4096.0334

$ws.Range("H34").Value = 3485.0244
$ws.Range("I34").Value = 3427.7273
$ws.Range("J34").Value = 3506.0334
$ws.Range("K34").Value = 3427.7273
$ws.Range("L34").Value = 3506.0334
$ws.Range("M34").Value = -3225.7273
$ws.Range("N34").Value = -3910.0334

$ws.Range("H99").Value = 3800
$ws.Range("I99").Value = 3050
$ws.Range("J99").Value = 5600
$ws.Range("K99").Value = 3050
$ws.Range("L99").Value = 5600
$ws.Range("M99").Value = -1552
$ws.Range("N99").Value = -8596

$ws.Range("H126").Value = 3800
$ws.Range("I126").Value = 3050
$ws.Range("J126").Value = 5600
$ws.Range("K126").Value = 9150
$ws.Range("L126").Value = 16800
$ws.Range("M126").Value = -6680
$ws.Range("N126").Value = -21740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 2366.6667
$ws.Range("I20").Value = 2366.6667
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 7100.000100000001
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -6873.000100000001
$ws.Range("N20").Value = ""

$ws.Range("H131").Value = 724.61
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 724.61
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2173.83
$ws.Range("M131").Value = ""
$ws.Range("N131").Value = -12253.83

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 15627209
$ws.Range("I102").Value = 18520572
$ws.Range("J102").Value = 3052.6
$ws.Range("K102").Value = 18520572
$ws.Range("L102").Value = 3052.6
$ws.Range("M102").Value = -18518950
$ws.Range("N102").Value = -6296.6

$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = -34940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 202.5
$ws.Range("I55").Value = 150
$ws.Range("J55").Value = 220
$ws.Range("K55").Value = 150
$ws.Range("L55").Value = 220
$ws.Range("M55").Value = 23
$ws.Range("N55").Value = -566

$ws.Range("H127").Value = 39852.824
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 39852.824
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 39852.824
$ws.Range("M127").Value = ""
$ws.Range("N127").Value = -49772.824

$ws.Range("H136").Value = 1367.2667
$ws.Range("I136").Value = 1393.5714
$ws.Range("J136").Value = 999
$ws.Range("K136").Value = 4180.7142
$ws.Range("L136").Value = 2997
$ws.Range("M136").Value = -1630.7142
$ws.Range("N136").Value = -8097

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1145.3572
$ws.Range("I126").Value = 1188.88
$ws.Range("J126").Value = 782.6667
$ws.Range("K126").Value = 3566.64
$ws.Range("L126").Value = 2348.0001
$ws.Range("M126").Value = -1096.64
$ws.Range("N126").Value = -7288.0001

$ws.Range("H132").Value = 1763.591
$ws.Range("I132").Value = 1061.6154
$ws.Range("J132").Value = 2777.5557
$ws.Range("K132").Value = 3184.8462
$ws.Range("L132").Value = 8332.667099999999
$ws.Range("M132").Value = -654.8462
$ws.Range("N132").Value = -13392.6671

$ws.Range("H136").Value = 21507728
$ws.Range("I136").Value = 31281654
$ws.Range("J136").Value = 5087.6665
$ws.Range("K136").Value = 93844962
$ws.Range("L136").Value = 15262.9995
$ws.Range("M136").Value = -93842412
$ws.Range("N136").Value = -20362.9995
